{"js": "// Statusbericht Woche 5 \u2014 apply two edits:\n//  1) Title line: merge \"Statusbericht \u2013 Woche 5\" + \", \" into a single run\n//     (text itself is unchanged: \"Statusbericht \u2013 Woche 5, 22.03.2022\").\n//  2) Body paragraph: replace the second half of the report text (starting\n//     with the old \"Wir sind bereit f\u00fcr den Sprint-Review...\" sentence)\n//     with the new Retrospective / Scrum-Poker text.\n\n// --- 1) Title paragraph -----------------------------------------------\nconst titleResults = context.document.body.search(\n  \"Statusbericht \u2013 Woche 5, \",\n  { matchCase: true }\n);\ntitleResults.load(\"items\");\nawait context.sync();\n\nif (titleResults.items.length > 0) {\n  const titleRange = titleResults.items[0];\n  // Re-insert the identical text so the run split collapses into one run,\n  // matching the merged run produced by the authored edit.\n  titleRange.insertText(\"Statusbericht \u2013 Woche 5, \", \"Replace\");\n  await context.sync();\n}\n\n// --- 2) Body paragraph ---------------------------------------------------\nconst oldTail =\n  \" Wir sind bereit f\u00fcr den Sprint-Review, welcher Heute Nachmittag am 22.03, ansteht um danach gut vorbereitet in den 2. Sprint zu starten.\";\nconst newTail =\n  \" Wir waren uns im Retrospective einig, dass unsere Herangehensweise an die Tasks grunds\u00e4tzlich gut war und wir diese f\u00fcr die 2. Iteration \u00fcbernehmen k\u00f6nnen. Zus\u00e4tzlich haben wir noch einige Ideen gesammelt\" +\n  \", welche unseren Arbeitsprozess verbessern k\u00f6nnten.\" +\n  \" Wir haben die Tasks f\u00fcr die 2. Iteration erhalten und deren Aufwand mit einem Scrum-Poker\" +\n  \" eingesch\u00e4tzt.\";\n\nconst bodyResults = context.document.body.search(oldTail, { matchCase: true });\nbodyResults.load(\"items\");\nawait context.sync();\n\nif (bodyResults.items.length > 0) {\n  const bodyRange = bodyResults.items[0];\n  bodyRange.insertText(newTail, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Statusbericht Woche 5 \u2014 apply two edits:\n#  1) Title line: merge \"Statusbericht \u2013 Woche 5\" + \", \" into a single run\n#     (visible text is unchanged: \"Statusbericht \u2013 Woche 5, 22.03.2022\").\n#  2) Body paragraph: replace the second half of the report text (starting\n#     with the old \"Wir sind bereit f\u00fcr den Sprint-Review...\" sentence)\n#     with the new Retrospective / Scrum-Poker text.\n\n$d = $word.ActiveDocument\n\n# --- 1) Title paragraph ---------------------------------------------------\n$titleRange = $d.Paragraphs(1).Range\n$titleFind = $titleRange.Find\n$titleFind.ClearFormatting()\n$titleFind.Text = \"Statusbericht \u2013 Woche 5, \"\n$titleFind.Replacement.ClearFormatting()\n$titleFind.Replacement.Text = \"Statusbericht \u2013 Woche 5, \"\n$titleFind.Execute($titleFind.Text, $false, $false, $false, $false, $false, $true, 1, $false, $titleFind.Replacement.Text, 1) | Out-Null\n\n# --- 2) Body paragraph -----------------------------------------------------\n$bodyRange = $d.Paragraphs(4).Range\n$bodyFind = $bodyRange.Find\n$bodyFind.ClearFormatting()\n$bodyFind.Text = \"Wir sind bereit f\u00fcr den Sprint-Review, welcher Heute Nachmittag am 22.03, ansteht um danach gut vorbereitet in den 2. Sprint zu starten.\"\n$bodyFind.Replacement.ClearFormatting()\n$bodyFind.Replacement.Text = \"Wir waren uns im Retrospective einig, dass unsere Herangehensweise an die Tasks grunds\u00e4tzlich gut war und wir diese f\u00fcr die 2. Iteration \u00fcbernehmen k\u00f6nnen. Zus\u00e4tzlich haben wir noch einige Ideen gesammelt, welche unseren Arbeitsprozess verbessern k\u00f6nnten. Wir haben die Tasks f\u00fcr die 2. Iteration erhalten und deren Aufwand mit einem Scrum-Poker eingesch\u00e4tzt.\"\n$bodyFind.Execute($bodyFind.Text, $false, $false, $false, $false, $false, $true, 1, $false, $bodyFind.Replacement.Text, 1) | Out-Null\n"}
